$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 39) with the next forecast entry, following the
# same layout as the existing data rows (date, year, y0 value, year, y1 value).
$ws.Range("A38").Copy()
$ws.Range("A39").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A39").Value = 45986

$ws.Range("B39").Value = 2025
$ws.Range("C39").Value = 0.2298740481777584
$ws.Range("D39").Value = 2026
$ws.Range("E39").Value = -0.05255865067609333
